# The commit inserts a new column ("Родитель" / Parent code) into the
# "Планы" sheet, right before the old column D ("Всего зачетных единиц"),
# shifting every subsequent column one place to the right. The new column
# holds the parent plan code (17916) for the three child plans (rows 3-5),
# while the root plan row (row 2) is left blank. Afterwards "Планы"
# becomes the active sheet with D5 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Планы")

# Insert a new blank column at D; everything from the old D onward shifts right.
$ws.Columns("D:D").Insert()

# New header cell.
$ws.Range("D1").Value = "Родитель"

# New parent-plan values (root row 2 stays empty).
$ws.Range("D3").Value = 17916
$ws.Range("D4").Value = 17916
$ws.Range("D5").Value = 17916

# Match the number formatting/style used by the sibling "Код" column (A).
$ws.Range("A3:A5").Copy()
$ws.Range("D3:D5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Give the new column its own width.
$ws.Columns("D:D").ColumnWidth = 14.1666666666667

# The author ended up viewing the "Планы" sheet with D5 selected.
$ws.Activate()
$ws.Range("D5").Select() | Out-Null
